$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update row 6 test data: email and password for TC001
$ws.Range("C6").Value = "sal@gmail.com"
$ws.Range("D6").Value = "Salsa123!"

# Update the active cell selection on the sheet
$ws.Range("G7").Select()

# Update the workbook window position
$excel.ActiveWindow.Left = 2730
$excel.ActiveWindow.Top = 2730
